$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Convert column C date-serial values (format "yyyy年m月") into literal text
# strings with the same displayed text, row by row, in sheet order. C3 already
# held the text "常住" and is re-entered in place so it lands right after the
# first converted date string in the shared-strings table, mirroring the diff.
$ws.Range("C2").Value = "2024年6月"
$ws.Range("C3").Value = "常住"
$ws.Range("C4").Value = "2023年7月"
$ws.Range("C5").Value = "2023年7月"
$ws.Range("C6").Value = "2024年4月"
$ws.Range("C7").Value = "2024年4月"
$ws.Range("C8").Value = "2024年3月"
$ws.Range("C9").Value = "2024年3月"
$ws.Range("C10").Value = "2024年4月"
$ws.Range("C11").Value = "2024年7月"
$ws.Range("C12").Value = "2024年7月"
$ws.Range("C13").Value = "2025年4月"
$ws.Range("C14").Value = "2025年4月"
$ws.Range("C15").Value = "2024年7月"
$ws.Range("C16").Value = "2024年6月"
$ws.Range("C17").Value = "2025年2月"
$ws.Range("C18").Value = "2023年1月"
$ws.Range("C19").Value = "2025年2月"
$ws.Range("C20").Value = "2025年2月"
$ws.Range("C21").Value = "2023年2月"
$ws.Range("C22").Value = "2025年2月"
$ws.Range("C23").Value = "2024年11月"
$ws.Range("C24").Value = "2024年7月"
$ws.Range("C25").Value = "2024年1月"
$ws.Range("C26").Value = "2023年9月"
$ws.Range("C27").Value = "2024年1月"
$ws.Range("C28").Value = "2023年7月"
$ws.Range("C29").Value = "2025年2月"
$ws.Range("C30").Value = "2024年10月"
$ws.Range("C31").Value = "2023年10月"
$ws.Range("C32").Value = "2023年10月"
$ws.Range("C33").Value = "2023年10月"
$ws.Range("C34").Value = "2025年1月"
$ws.Range("C35").Value = "2025年1月"
$ws.Range("C36").Value = "2024年10月"
$ws.Range("C37").Value = "2023年7月"
$ws.Range("C38").Value = "2023年7月"
$ws.Range("C39").Value = "2023年7月"
$ws.Range("C40").Value = "2025年1月"
$ws.Range("C41").Value = "2025年1月"
$ws.Range("C42").Value = "2023年10月"
$ws.Range("C43").Value = "2024年1月"
$ws.Range("C44").Value = "2025年6月"
$ws.Range("C45").Value = "2023年4月"
$ws.Range("C46").Value = "2024年10月"
$ws.Range("C47").Value = "2024年10月"
$ws.Range("C48").Value = "2025年1月"
$ws.Range("C49").Value = "2024年3月"
$ws.Range("C50").Value = "2023年4月"
$ws.Range("C51").Value = "2023年8月"
$ws.Range("C52").Value = "2025年2月"
$ws.Range("C53").Value = "2024年9月"
$ws.Range("C54").Value = "2023年10月"
$ws.Range("C55").Value = "2024年10月"
$ws.Range("C56").Value = "2024年3月"
$ws.Range("C57").Value = "2023年6月"
$ws.Range("C58").Value = "2024年7月"
$ws.Range("C59").Value = "2024年10月"
$ws.Range("C60").Value = "2024年3月"
$ws.Range("C61").Value = "2024年3月"
$ws.Range("C62").Value = "2024年4月"
$ws.Range("C63").Value = "2024年1月"
$ws.Range("C64").Value = "2024年1月"
$ws.Range("C65").Value = "2023年6月"
$ws.Range("C66").Value = "2023年6月"
$ws.Range("C67").Value = "2023年9月"
$ws.Range("C68").Value = "2024年8月"
$ws.Range("C69").Value = "2024年1月"
$ws.Range("C70").Value = "2024年8月"
$ws.Range("C71").Value = "2024年8月"
$ws.Range("C72").Value = "2024年1月"
$ws.Range("C73").Value = "2024年7月"
$ws.Range("C74").Value = "2023年7月"
$ws.Range("C75").Value = "2025年3月"
$ws.Range("C76").Value = "2024年1月"
$ws.Range("C77").Value = "2024年7月"
$ws.Range("C78").Value = "2024年1月"
$ws.Range("C79").Value = "2024年7月"
$ws.Range("C80").Value = "2022年12月"
$ws.Range("C81").Value = "2024年7月"
$ws.Range("C82").Value = "2024年7月"
$ws.Range("C83").Value = "2024年7月"
$ws.Range("C84").Value = "2024年8月"
$ws.Range("C85").Value = "2024年8月"
$ws.Range("C86").Value = "2024年8月"
$ws.Range("C87").Value = "2025年4月"
$ws.Range("C88").Value = "2025年2月"
$ws.Range("C89").Value = "2025年2月"
$ws.Range("C90").Value = "2024年1月"
$ws.Range("C91").Value = "2024年1月"
$ws.Range("C92").Value = "2024年4月"
$ws.Range("C93").Value = "2025年2月"
$ws.Range("C94").Value = "2025年2月"
$ws.Range("C95").Value = "2025年2月"
$ws.Range("C96").Value = "2025年2月"
$ws.Range("C97").Value = "2024年1月"
$ws.Range("C98").Value = "2024年4月"
$ws.Range("C99").Value = "2024年1月"
$ws.Range("C100").Value = "2024年1月"
$ws.Range("C101").Value = "2024年1月"
$ws.Range("C102").Value = "2025年7月"
$ws.Range("C103").Value = "2025年7月"

# Column C was manually resized (drag) instead of using "best fit"
$ws.Columns("C").ColumnWidth = 10.6

# Selection moved from B1:B1048576 (whole column) to C3, and the frozen/
# scrolled top-left cell reset back to A1
$ws.Range("A1").Select()
$ws.Range("C3").Select()
